$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Every existing data row (2-429) has its "Förändrad" (column C) date
#    bumped from 2023-09-13 (45182) to 2023-09-15 (45184).
for ($r = 2; $r -le 429; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}

# 2) Row 429 previously lacked an explicit row height; give it the same
#    15pt custom height used by all the other data rows.
$ws.Rows.Item(429).RowHeight = 15

# 3) Append two brand-new felling notifications as rows 430 and 431.

# Row 430: A 42806-2023
$ws.Cells.Item(430, 1).Value = "A 42806-2023"
$ws.Cells.Item(430, 2).Value = 45182
$ws.Cells.Item(430, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(430, 3).Value = 45184
$ws.Cells.Item(430, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(430, 4).Value = "SÖDERMANLANDS LÄN"
$ws.Cells.Item(430, 5).Value = "STRÄNGNÄS"
$ws.Cells.Item(430, 7).Value = 4.7
$ws.Cells.Item(430, 8).Value = 0
$ws.Cells.Item(430, 9).Value = 0
$ws.Cells.Item(430, 10).Value = 0
$ws.Cells.Item(430, 11).Value = 0
$ws.Cells.Item(430, 12).Value = 0
$ws.Cells.Item(430, 13).Value = 0
$ws.Cells.Item(430, 14).Value = 0
$ws.Cells.Item(430, 15).Value = 0
$ws.Cells.Item(430, 16).Value = 0
$ws.Cells.Item(430, 17).Value = 0
$ws.Cells.Item(430, 18).WrapText = $true
$ws.Rows.Item(430).RowHeight = 15

# Row 431: A 42808-2023
$ws.Cells.Item(431, 1).Value = "A 42808-2023"
$ws.Cells.Item(431, 2).Value = 45182
$ws.Cells.Item(431, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(431, 3).Value = 45184
$ws.Cells.Item(431, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(431, 4).Value = "SÖDERMANLANDS LÄN"
$ws.Cells.Item(431, 5).Value = "STRÄNGNÄS"
$ws.Cells.Item(431, 7).Value = 1.4
$ws.Cells.Item(431, 8).Value = 0
$ws.Cells.Item(431, 9).Value = 0
$ws.Cells.Item(431, 10).Value = 0
$ws.Cells.Item(431, 11).Value = 0
$ws.Cells.Item(431, 12).Value = 0
$ws.Cells.Item(431, 13).Value = 0
$ws.Cells.Item(431, 14).Value = 0
$ws.Cells.Item(431, 15).Value = 0
$ws.Cells.Item(431, 16).Value = 0
$ws.Cells.Item(431, 17).Value = 0
$ws.Cells.Item(431, 18).WrapText = $true
